# Swap the taxon data between row 5 and row 6.
# Columns involved: A (Id), B (Taxonsorteringsordning), D (Rodlistade),
# E (TaxonId), F (Artnamn), G (Vetenskapligt namn), H (Auktor),
# Q (Ost), R (Nord).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $cell5 = $ws.Range($col + "5")
    $cell6 = $ws.Range($col + "6")

    $v5 = $cell5.Value2
    $v6 = $cell6.Value2

    $cell5.Value2 = $v6
    $cell6.Value2 = $v5
}
